# Removed an extra slide (the last slide, a blank "Title 1"-only slide)
# from the presentation.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item($p.Slides.Count)
$s.Delete()
